$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 625.9259
$ws.Range("J17").Value = 625.9259
$ws.Range("L17").Value = 1877.7777
$ws.Range("N17").Value = -2213.7777

$ws.Range("H33").Value = 169.90909
$ws.Range("I33").Value = 46.9
$ws.Range("J33").Value = 1400
$ws.Range("K33").Value = 46.9
$ws.Range("L33").Value = 1400
$ws.Range("M33").Value = 182.1
$ws.Range("N33").Value = -1858

$ws.Range("H40").Value = 1816.3334
$ws.Range("I40").Value = 1787.625
$ws.Range("K40").Value = 1787.625
$ws.Range("M40").Value = -1612.625

$ws.Range("H76").Value = 3862.125
$ws.Range("I76").Value = 4000.5
$ws.Range("J76").Value = 3816
$ws.Range("K76").Value = 4000.5
$ws.Range("L76").Value = 3816
$ws.Range("M76").Value = -3685.5
$ws.Range("N76").Value = -4446

$ws.Range("H79").Value = 3862.125
$ws.Range("I79").Value = 4000.5
$ws.Range("J79").Value = 3816
$ws.Range("K79").Value = 4000.5
$ws.Range("L79").Value = 3816
$ws.Range("M79").Value = -2908.5
$ws.Range("N79").Value = -6000

$ws.Range("H100").Value = 1409.125
$ws.Range("I100").Value = 1444.6
$ws.Range("J100").Value = 1350
$ws.Range("K100").Value = 1444.6
$ws.Range("L100").Value = 1350
$ws.Range("M100").Value = -903.5999999999999
$ws.Range("N100").Value = -2432

$ws.Range("H106").Value = 4026.4167
$ws.Range("I106").Value = 2386.4285
$ws.Range("J106").Value = 6322.4
$ws.Range("K106").Value = 2386.4285
$ws.Range("L106").Value = 6322.4
$ws.Range("M106").Value = -1755.4285
$ws.Range("N106").Value = -7584.4

$ws.Range("H125").Value = 4950
$ws.Range("J125").Value = 4950
$ws.Range("L125").Value = 44550
$ws.Range("N125").Value = -49470

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 50000250
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 100000000
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 100000000
$ws.Range("M3").Value = -385
$ws.Range("N3").Value = -100000230

$ws.Range("H45").Value = 1610.7106
$ws.Range("I45").Value = 1026
$ws.Range("J45").Value = 1742.742
$ws.Range("K45").Value = 1026
$ws.Range("L45").Value = 1742.742
$ws.Range("M45").Value = -649
$ws.Range("N45").Value = -2496.742

$ws.Range("H101").Value = 36500
$ws.Range("J101").Value = 36500
$ws.Range("L101").Value = 36500
$ws.Range("N101").Value = -42990

$ws.Range("H102").Value = 1118.5714
$ws.Range("I102").Value = 1034.5454
$ws.Range("J102").Value = 1426.6666
$ws.Range("K102").Value = 1034.5454
$ws.Range("L102").Value = 1426.6666
$ws.Range("M102").Value = 587.4546
$ws.Range("N102").Value = -4670.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 800
$ws.Range("I8").Value = 510
$ws.Range("J8").Value = 2250
$ws.Range("K8").Value = 510
$ws.Range("L8").Value = 2250
$ws.Range("M8").Value = -370
$ws.Range("N8").Value = -2530

$ws.Range("H94").Value = 432.05554
$ws.Range("I94").Value = 467.64285
$ws.Range("J94").Value = 307.5
$ws.Range("K94").Value = 467.64285
$ws.Range("L94").Value = 307.5
$ws.Range("M94").Value = -16.64285000000001
$ws.Range("N94").Value = -1209.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1822.2766
$ws.Range("I31").Value = 1427.625
$ws.Range("J31").Value = 4077.4285
$ws.Range("K31").Value = 1427.625
$ws.Range("L31").Value = 4077.4285
$ws.Range("M31").Value = -1132.625
$ws.Range("N31").Value = -4667.4285

$ws.Range("H34").Value = 1822.2766
$ws.Range("I34").Value = 1427.625
$ws.Range("J34").Value = 4077.4285
$ws.Range("K34").Value = 1427.625
$ws.Range("L34").Value = 4077.4285
$ws.Range("M34").Value = -1225.625
$ws.Range("N34").Value = -4481.4285

$ws.Range("H132").Value = 2195.12
$ws.Range("I132").Value = 1675.7273
$ws.Range("J132").Value = 6004
$ws.Range("K132").Value = 5027.1819
$ws.Range("L132").Value = 18012
$ws.Range("M132").Value = -2497.1819
$ws.Range("N132").Value = -23072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 3529.7144
$ws.Range("I115").Value = 528
$ws.Range("J115").Value = 4030
$ws.Range("K115").Value = 1584
$ws.Range("L115").Value = 12090
$ws.Range("M115").Value = -409
$ws.Range("N115").Value = -14440

$ws.Range("H116").Value = 2458.9
$ws.Range("I116").Value = 1226.8572
$ws.Range("J116").Value = 5333.6665
$ws.Range("K116").Value = 3680.5716
$ws.Range("L116").Value = 16000.9995
$ws.Range("M116").Value = -238.5715999999998
$ws.Range("N116").Value = -22884.9995

$ws.Range("H121").Value = 708.2353000000001
$ws.Range("I121").Value = 440
$ws.Range("J121").Value = 5000
$ws.Range("K121").Value = 1320
$ws.Range("L121").Value = 15000
$ws.Range("M121").Value = -10
$ws.Range("N121").Value = -17620

$ws.Range("H129").Value = 5564.769
$ws.Range("J129").Value = 9176.538
$ws.Range("L129").Value = 27529.614
$ws.Range("N129").Value = -37529.614

$ws.Range("H131").Value = 856.36664
$ws.Range("I131").Value = 684
$ws.Range("J131").Value = 890.84
$ws.Range("K131").Value = 2052
$ws.Range("L131").Value = 2672.52
$ws.Range("M131").Value = 2988
$ws.Range("N131").Value = -12752.52

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1511.8
$ws.Range("I97").Value = 1511.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1511.8
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1015.8
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 2515.4
$ws.Range("I102").Value = 2440.8235
$ws.Range("K102").Value = 2440.8235
$ws.Range("M102").Value = -818.8235

$ws.Range("H122").Value = 9092451
$ws.Range("I122").Value = 11112329
$ws.Range("K122").Value = 33336987
$ws.Range("M122").Value = -33334537

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4103.8335
$ws.Range("I40").Value = 4128.4287
$ws.Range("J40").Value = 4069.4
$ws.Range("K40").Value = 4128.4287
$ws.Range("L40").Value = 4069.4
$ws.Range("M40").Value = -3992.4287
$ws.Range("N40").Value = -4341.4

$ws.Range("H104").Value = 12296.2
$ws.Range("J104").Value = 12296.2
$ws.Range("L104").Value = 12296.2
$ws.Range("N104").Value = -19284.2

$ws.Range("H122").Value = 2326
$ws.Range("I122").Value = 1998
$ws.Range("J122").Value = 2490
$ws.Range("K122").Value = 5994
$ws.Range("L122").Value = 7470
$ws.Range("M122").Value = -3544
$ws.Range("N122").Value = -12370

$ws.Range("H132").Value = 3392.5588
$ws.Range("I132").Value = 2807.125
$ws.Range("J132").Value = 4797.6
$ws.Range("K132").Value = 8421.375
$ws.Range("L132").Value = 14392.8
$ws.Range("M132").Value = -5891.375
$ws.Range("N132").Value = -19452.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 19938.5
$ws.Range("J54").Value = 19938.5
$ws.Range("L54").Value = 19938.5
$ws.Range("N54").Value = -20978.5

$ws.Range("H81").Value = 2371.5715
$ws.Range("I81").Value = 2250.25
$ws.Range("J81").Value = 2533.3333
$ws.Range("K81").Value = 4500.5
$ws.Range("L81").Value = 5066.6666
$ws.Range("M81").Value = -3439.5
$ws.Range("N81").Value = -7188.6666

$ws.Range("H84").Value = 2371.5715
$ws.Range("I84").Value = 2250.25
$ws.Range("J84").Value = 2533.3333
$ws.Range("K84").Value = 22502.5
$ws.Range("L84").Value = 25333.333
$ws.Range("M84").Value = -17198.5
$ws.Range("N84").Value = -35941.333

$ws.Range("H96").Value = 2566.9333
$ws.Range("I96").Value = 2482.8572
$ws.Range("J96").Value = 2640.5
$ws.Range("K96").Value = 2482.8572
$ws.Range("L96").Value = 2640.5
$ws.Range("M96").Value = -1109.8572
$ws.Range("N96").Value = -5386.5
